# Generate Report for Handoff
# - Flip the localization status text from "In Translation" to
#   "Ready for handoff" everywhere it appears (Overview!E2:F2,
#   zh-cn!C2, de-de!C2 all share the same text).
# - Bump the two "Latest * Datetime" timestamps that recorded this event.
# - Widen the Status/zh-cn/de-de status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status text: "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# Timestamps
$overview.Range("G2").Value = "2016-08-31 09:14:10"
$zhcn.Range("H2").Value     = "2016-08-31 09:13:59"
$dede.Range("H2").Value     = "2016-08-31 09:14:10"

# Widen the status columns so the longer "Ready for handoff" text fits.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33
